# Daily attendance processing - swap the order of "Recorded By" entries
# from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# in column G ("Recorded By") of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2
    if ($value -eq $oldText) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Output "Updated $changed 'Recorded By' cells in column G."
